$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.635.25'
$ws.Range("E2").Value = '  +0.69%  '
$ws.Range("D3").Value = '1.844.54'
$ws.Range("E3").Value = '  +0.05%  '
$ws.Range("E4").Value = '  +0.09%  '
$ws.Range("E5").Value = '  -0.78%  '
$ws.Range("E6").Value = '  +0.06%  '
$ws.Range("D7").Value = '0.5271'
$ws.Range("E7").Value = '  +1.18%  '
$ws.Range("D8").Value = '0.3156'
$ws.Range("E8").Value = '  -3.45%  '
$ws.Range("D9").Value = '0.06808'
$ws.Range("E9").Value = '  +0.16%  '
$ws.Range("D10").Value = '19.13'
$ws.Range("E10").Value = '  +2.10%  '
$ws.Range("D11").Value = '0.7868'
$ws.Range("E11").Value = '  +0.87%  '
$ws.Range("D12").Value = '0.07783'
$ws.Range("E12").Value = '  +0.38%  '
$ws.Range("D13").Value = '1.837.53'
$ws.Range("E13").Value = '  -0.45%  '
$ws.Range("D14").Value = '88.44'
$ws.Range("E14").Value = '  +0.52%  '
$ws.Range("E15").Value = '  +0.08%  '
$ws.Range("D16").Value = '1.000'
$ws.Range("E16").Value = '  +0.12%  '
$ws.Range("E17").Value = '  -0.01%  '
$ws.Range("E18").Value = '  +0.08%  '
$ws.Range("D19").Value = '0.000007927'
$ws.Range("E19").Value = '  -0.61%  '
$ws.Range("D20").Value = '26.667.17'
$ws.Range("E20").Value = '  +0.69%  '
$ws.Range("D21").Value = '2.080.29'
$ws.Range("E21").Value = '  +0.35%  '
$ws.Range("D22").Value = '4.616'
$ws.Range("E22").Value = '  -0.11%  '
$ws.Range("D23").Value = '6.005'
$ws.Range("E23").Value = '  +0.41%  '
$ws.Range("D24").Value = '9.339'
$ws.Range("E24").Value = '  -2.27%  '
$ws.Range("B25").Value = 'Monero'
$ws.Range("C25").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D25").Value = '143.17'
$ws.Range("E25").Value = '  -0.99%  '
$ws.Range("B26").Value = 'LidoDAOToken'
$ws.Range("C26").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D26").Value = '2.225'
$ws.Range("E26").Value = '  +2.08%  '
$ws.Range("E27").Value = '  +2.28%  '
$ws.Range("D28").Value = '17.06'
$ws.Range("E28").Value = '  +0.38%  '
$ws.Range("D29").Value = '111.05'
$ws.Range("E29").Value = '  -0.77%  '
$ws.Range("D30").Value = '4.220'
$ws.Range("E30").Value = '  +1.50%  '
$ws.Range("D31").Value = '0.08706'
$ws.Range("E31").Value = '  -0.01%  '
$ws.Range("D32").Value = '4.090'
$ws.Range("E32").Value = '  -0.92%  '
$ws.Range("D33").Value = '0.04886'
$ws.Range("E33").Value = '  +0.98%  '
$ws.Range("D34").Value = '0.7312'
$ws.Range("E34").Value = '  +1.01%  '
$ws.Range("D35").Value = '1.140'
$ws.Range("E35").Value = '  +0.78%  '
$ws.Range("E36").Value = '  +0.77%  '
$ws.Range("D37").Value = '3.116'
$ws.Range("E37").Value = '  +0.07%  '
$ws.Range("D38").Value = '2.363'
$ws.Range("E38").Value = '  +5.84%  '
$ws.Range("E39").Value = '  -2.32%  '
$ws.Range("D40").Value = '0.4836'
$ws.Range("E40").Value = '  -0.22%  '
$ws.Range("E41").Value = '  -0.23%  '
$ws.Range("D42").Value = '109.68'
$ws.Range("E42").Value = '  -1.37%  '
$ws.Range("D43").Value = '5.934'
$ws.Range("E43").Value = '  -2.28%  '
$ws.Range("D44").Value = '1.000'
$ws.Range("E44").Value = '  +0.10%  '
$ws.Range("D45").Value = '7.743'
$ws.Range("E45").Value = '  -0.14%  '
$ws.Range("D46").Value = '0.4211'
$ws.Range("E46").Value = '  +0.89%  '
$ws.Range("D47").Value = '9.055'
$ws.Range("E47").Value = '  -0.15%  '
$ws.Range("D48").Value = '0.1245'
$ws.Range("E48").Value = '  +0.91%  '
$ws.Range("B49").Value = 'Elrond'
$ws.Range("C49").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D49").Value = '34.94'
$ws.Range("E49").Value = '  -0.48%  '
$ws.Range("B50").Value = 'Cronos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D50").Value = '0.05828'
$ws.Range("E50").Value = '  -1.73%  '
$ws.Range("D51").Value = '0.8983'
$ws.Range("E51").Value = '  +1.01%  '
